$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-28 Sunday", "2024-07-29 Monday"),
    @("98×53=", "28×17="),
    @("28×89=", "16×28="),
    @("11×21=", "98×79="),
    @("67×38=", "18×77="),
    @("35×86=", "29×32="),
    @("15×89=", "87×98="),
    @("96×28=", "49×53="),
    @("68×25=", "56×73="),
    @("19×85=", "94×12="),
    @("20×26=", "28×84="),
    @("29×36=", "60×39="),
    @("21×35=", "41×37="),
    @("86×78=", "30×51="),
    @("14×97=", "63×43="),
    @("98×55=", "43×84="),
    @("74×29=", "27×69="),
    @("59×52=", "46×78="),
    @("59×13=", "32×97="),
    @("64×18=", "23×69="),
    @("53×44=", "50×45="),
    @("27×83=", "56×82="),
    @("28×57=", "30×40="),
    @("23×25=", "54×23="),
    @("44×78=", "34×40="),
    @("45×98=", "17×71=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
